$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("150:150").Insert()

$ws.Range("A150").Value = 10
$ws.Range("B150").Value = "Vega Modelo de Temuco"
$ws.Range("C150").Value = "La Araucanía"
$ws.Range("D150").Value = 45089
$ws.Range("E150").Value = 9
$ws.Range("F150").Value = "Fruta"
$ws.Range("G150").Value = 100104
$ws.Range("H150").Value = "Frutos de pepita"
$ws.Range("I150").Value = 100104001
$ws.Range("J150").Value = "Granada"
$ws.Range("K150").Value = "Wonderfull"
$ws.Range("L150").Value = "Primera"
$ws.Range("M150").Value = 500
$ws.Range("N150").Value = 12000
$ws.Range("O150").Value = 13000
$ws.Range("P150").Value = 12400
$ws.Range("Q150").Value = "$/bandeja 10 kilos granel"
$ws.Range("R150").Value = "Provincia de Limarí"
$ws.Range("S150").Value = 1240
$ws.Range("T150").Value = 10
